$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 14 de Octubre de 2020 a las 02:13"

# Row 4
$ws.Range("B4").Value = 8087865
$ws.Range("C4").Value = 49289
$ws.Range("D4").Value = 5221056
$ws.Range("E4").Value = 2646007
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 784
$ws.Range("H4").Value = 220802

# Row 6
$ws.Range("B6").Value = 5114823
$ws.Range("C6").Value = 11415
$ws.Range("D6").Value = 4526975
$ws.Range("E6").Value = 436785
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 354
$ws.Range("H6").Value = 151063

# Row 10
$ws.Range("B10").Value = 917035
$ws.Range("C10").Value = 13305
$ws.Range("D10").Value = 742235
$ws.Range("E10").Value = 150228
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 386
$ws.Range("H10").Value = 24572

# Row 11
$ws.Range("B11").Value = 853974
$ws.Range("C11").Value = 2803
$ws.Range("D11").Value = 753959
$ws.Range("E11").Value = 66596
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 62
$ws.Range("H11").Value = 33419

# Row 30
$ws.Range("B30").Value = 186881
$ws.Range("C30").Value = 4042
$ws.Range("D30").Value = 157486
$ws.Range("E30").Value = 19741
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 27
$ws.Range("H30").Value = 9654

# Row 31
$ws.Range("D31").Value = 20325
$ws.Range("E31").Value = 135344

# Row 37
$ws.Range("A37").Value = "Chequia"
$ws.Range("B37").Value = 129747
$ws.Range("C37").Value = 8326
$ws.Range("D37").Value = 59901
$ws.Range("E37").Value = 68740
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 55
$ws.Range("H37").Value = 1106

# Row 38
$ws.Range("A38").Value = "Catar"
$ws.Range("B38").Value = 128405
$ws.Range("C38").Value = 214
$ws.Range("D38").Value = 125373
$ws.Range("E38").Value = 2812
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 220

# Row 39
$ws.Range("B39").Value = 121296
$ws.Range("C39").Value = 494
$ws.Range("D39").Value = 96675
$ws.Range("E39").Value = 22110
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 9
$ws.Range("H39").Value = 2511

# Row 67
$ws.Range("B67").Value = 51197
$ws.Range("C67").Value = 853
$ws.Range("D67").Value = 33325
$ws.Range("E67").Value = 16764
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 12
$ws.Range("H67").Value = 1108

# Row 69
$ws.Range("B69").Value = 47126
$ws.Range("C69").Value = 96
$ws.Range("D69").Value = 46469
$ws.Range("E69").Value = 347
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 310

# Row 131
$ws.Range("A131").Value = "Surinam"
$ws.Range("B131").Value = 5072
$ws.Range("C131").Value = 14
$ws.Range("D131").Value = 4870
$ws.Range("E131").Value = 95
$ws.Range("H131").Value = 107

# Row 132
$ws.Range("A132").Value = "Guinea Ecuatorial"
$ws.Range("B132").Value = 5066
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 4954
$ws.Range("E132").Value = 29
$ws.Range("H132").Value = 83

# Row 156
$ws.Range("B156").Value = 2337
$ws.Range("C156").Value = 24
$ws.Range("D156").Value = 1987
$ws.Range("E156").Value = 299

# Row 159
$ws.Range("A159").Value = "Martinica"
$ws.Range("B159").Value = 2257
$ws.Range("C159").Value = 406
$ws.Range("D159").Value = 98
$ws.Range("E159").Value = 2135
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 2
$ws.Range("H159").Value = 24

# Row 160
$ws.Range("A160").Value = "Republica de Chipre"
$ws.Range("B160").Value = 2130
$ws.Range("C160").Value = 83
$ws.Range("D160").Value = 1444
$ws.Range("E160").Value = 661
$ws.Range("H160").Value = 25

# Row 161
$ws.Range("A161").Value = "Yemen"
$ws.Range("B161").Value = 2053
$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 1329
$ws.Range("E161").Value = 128
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 596

# Row 162
$ws.Range("A162").Value = "Togo"
$ws.Range("B162").Value = 1972
$ws.Range("C162").Value = 23
$ws.Range("D162").Value = 1465
$ws.Range("E162").Value = 457
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 1
$ws.Range("H162").Value = 50

# Row 163
$ws.Range("A163").Value = "Nueva Zelanda"
$ws.Range("B163").Value = 1872
$ws.Range("C163").Value = 1
$ws.Range("D163").Value = 1808
$ws.Range("E163").Value = 39
$ws.Range("H163").Value = 25

# Row 171
$ws.Range("A171").Value = "San Martin (Parte Holandesa)"
$ws.Range("B171").Value = 719
$ws.Range("C171").Value = 9
$ws.Range("D171").Value = 645
$ws.Range("E171").Value = 52
$ws.Range("H171").Value = 22

# Row 172
$ws.Range("A172").Value = "Crucero"
$ws.Range("B172").Value = 712
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 659
$ws.Range("E172").Value = 40
$ws.Range("H172").Value = 13

# Row 173
$ws.Range("D173").Value = 673
$ws.Range("E173").Value = 17

# Row 192
$ws.Range("B192").Value = 185
$ws.Range("C192").Value = 1
$ws.Range("D192").Value = 172
$ws.Range("E192").Value = 4
